# "SOP. CAÑOS CURVOS" price-list sheet (Hoja1):
#  - bump the date serial in A1 by one day (45310 -> 45311)
#  - update the two unit prices in column D (rows 30 and 31)
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

$ws.Range("A1").Value = 45311
$ws.Range("D30").Value = 936
$ws.Range("D31").Value = 1040
